{"js": "// Replace each two-digit-divided-by-one-digit equation's source text\n// with the regenerated equation at the same table position, in\n// document order (this mirrors the canonical OOXML diff exactly,\n// since none of the 25 old equation strings repeat elsewhere in the body).\nconst replacements = [\n  [\"71\u00f79=\", \"23\u00f76=\"],\n  [\"16\u00f76=\", \"28\u00f75=\"],\n  [\"98\u00f76=\", \"13\u00f74=\"],\n  [\"65\u00f78=\", \"93\u00f73=\"],\n  [\"13\u00f75=\", \"46\u00f74=\"],\n  [\"41\u00f78=\", \"18\u00f75=\"],\n  [\"64\u00f75=\", \"58\u00f73=\"],\n  [\"93\u00f75=\", \"80\u00f78=\"],\n  [\"95\u00f76=\", \"92\u00f73=\"],\n  [\"59\u00f75=\", \"32\u00f78=\"],\n  [\"16\u00f75=\", \"74\u00f77=\"],\n  [\"74\u00f75=\", \"75\u00f76=\"],\n  [\"50\u00f78=\", \"86\u00f78=\"],\n  [\"65\u00f76=\", \"63\u00f73=\"],\n  [\"63\u00f75=\", \"75\u00f74=\"],\n  [\"63\u00f76=\", \"52\u00f78=\"],\n  [\"15\u00f75=\", \"66\u00f77=\"],\n  [\"44\u00f73=\", \"11\u00f72=\"],\n  [\"63\u00f79=\", \"78\u00f74=\"],\n  [\"31\u00f78=\", \"83\u00f73=\"],\n  [\"89\u00f78=\", \"32\u00f79=\"],\n  [\"63\u00f77=\", \"24\u00f74=\"],\n  [\"19\u00f74=\", \"61\u00f75=\"],\n  [\"47\u00f75=\", \"64\u00f77=\"],\n  [\"94\u00f74=\", \"94\u00f73=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  // eslint-disable-next-line no-await-in-loop\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`No match found for \"${oldText}\"`);\n  }\n\n  // Each equation string is unique in the document, so there should\n  // be exactly one hit; replace it (or all hits, defensively) in place.\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace each two-digit-divided-by-one-digit equation's source text\n# with the regenerated equation at the same table position, in\n# document order. Each of the 25 old equation strings is unique in\n# the document body, so a simple Find/Replace per pair is safe and\n# mirrors the canonical OOXML diff exactly.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"71\u00f79=\", \"23\u00f76=\"),\n  @(\"16\u00f76=\", \"28\u00f75=\"),\n  @(\"98\u00f76=\", \"13\u00f74=\"),\n  @(\"65\u00f78=\", \"93\u00f73=\"),\n  @(\"13\u00f75=\", \"46\u00f74=\"),\n  @(\"41\u00f78=\", \"18\u00f75=\"),\n  @(\"64\u00f75=\", \"58\u00f73=\"),\n  @(\"93\u00f75=\", \"80\u00f78=\"),\n  @(\"95\u00f76=\", \"92\u00f73=\"),\n  @(\"59\u00f75=\", \"32\u00f78=\"),\n  @(\"16\u00f75=\", \"74\u00f77=\"),\n  @(\"74\u00f75=\", \"75\u00f76=\"),\n  @(\"50\u00f78=\", \"86\u00f78=\"),\n  @(\"65\u00f76=\", \"63\u00f73=\"),\n  @(\"63\u00f75=\", \"75\u00f74=\"),\n  @(\"63\u00f76=\", \"52\u00f78=\"),\n  @(\"15\u00f75=\", \"66\u00f77=\"),\n  @(\"44\u00f73=\", \"11\u00f72=\"),\n  @(\"63\u00f79=\", \"78\u00f74=\"),\n  @(\"31\u00f78=\", \"83\u00f73=\"),\n  @(\"89\u00f78=\", \"32\u00f79=\"),\n  @(\"63\u00f77=\", \"24\u00f74=\"),\n  @(\"19\u00f74=\", \"61\u00f75=\"),\n  @(\"47\u00f75=\", \"64\u00f77=\"),\n  @(\"94\u00f74=\", \"94\u00f73=\"),\n)\n\nforeach ($pair in $pairs) {\n  $old = $pair[0]\n  $new = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $old\n  $find.Replacement.Text = $new\n  $find.Forward = $true\n  $find.Wrap = 1           # wdFindContinue\n  $found = $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)  # wdReplaceAll\n\n  if (-not $found) {\n    throw \"No match found for $old\"\n  }\n}\n"}
